# The workbook contains a daily price log for "Vega Modelo de Temuco - Zanahoria".
# A new daily record needs to be inserted as row 478, pushing the former rows
# 478-523 down to 479-524 (the former last row ends up as the new row 524).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 478, shifting existing rows down.
$ws.Rows.Item(478).Insert()

# Populate the new row with the new record's data.
$ws.Cells.Item(478, 1).Value = 10
$ws.Cells.Item(478, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(478, 3).Value = "La Araucanía"
$ws.Cells.Item(478, 4).Value = 45166
$ws.Cells.Item(478, 5).Value = 9
$ws.Cells.Item(478, 6).Value = 100114013
$ws.Cells.Item(478, 7).Value = "Zanahoria"
$ws.Cells.Item(478, 8).Value = "Sin especificar"
$ws.Cells.Item(478, 9).Value = "Primera"
$ws.Cells.Item(478, 10).Value = 200
$ws.Cells.Item(478, 11).Value = 5000
$ws.Cells.Item(478, 12).Value = 5000
$ws.Cells.Item(478, 13).Value = 5000
$ws.Cells.Item(478, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(478, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(478, 16).Value = 200
$ws.Cells.Item(478, 17).Value = 25
$ws.Cells.Item(478, 18).Value = "Hortaliza"

# Make sure the date cell keeps the existing date number format used by the
# rest of column D (style carried over automatically from the row insert,
# but set explicitly to be safe).
$ws.Cells.Item(478, 4).NumberFormat = $ws.Cells.Item(479, 4).NumberFormat
